$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (names / card number) ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit value that must remain TEXT (it was stored as text
# before the edit). A leading apostrophe forces Excel to keep it as text
# instead of silently converting it to a Number. That, however, marks the
# cell with a "quote prefix" style flag, which would bump its style index.
# Restore the original formatting by pasting just the format from a
# neighbouring cell that already shares B3's original style.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 11.09.2023"

# --- Transaction rows 6-11 ---
$ws.Range("B6").Value = "13.09."
$ws.Range("C6").Value = "14.09."
$ws.Range("D6").Value = "PAYPAL XGQIZK"
$ws.Range("E6").Value = "30,86-"

$ws.Range("B7").Value = "16.09."
$ws.Range("C7").Value = "17.09."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-94465518"
$ws.Range("E7").Value = "56,42-"

$ws.Range("B8").Value = "19.09."
$ws.Range("C8").Value = "20.09."
$ws.Range("D8").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E8").Value = "47,09-"

$ws.Range("B9").Value = "22.09."
$ws.Range("C9").Value = "23.09."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 12249199"
$ws.Range("E9").Value = "84,59-"

$ws.Range("B10").Value = "24.09."
$ws.Range("C10").Value = "25.09."
$ws.Range("D10").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E10").Value = "82,37-"

$ws.Range("B11").Value = "27.09."
$ws.Range("C11").Value = "28.09."
$ws.Range("D11").Value = "KARTENZ./27.09 EDEKA RO"
$ws.Range("E11").Value = "116,02-"

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 01.10.2023"
$ws.Range("E12").Value = "417,35-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 10.10.2023"
